$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '69.764.84'
$ws.Range("E2").Value = '  +6.02%  '
$ws.Range("D3").Value = '3.576.75'
$ws.Range("E3").Value = '  +5.07%  '
$ws.Range("E4").Value = '  +0.11%  '
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = '592.21'
$c.ClearFormats()
$ws.Range("E5").Value = '  +5.57%  '
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = '192.17'
$c.ClearFormats()
$ws.Range("E6").Value = '  +8.89%  '
$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = '0.641'
$c.ClearFormats()
$ws.Range("E7").Value = '  +1.46%  '
$ws.Range("D8").Value = '3.569.20'
$ws.Range("E8").Value = '  +5.09%  '
$ws.Range("E9").Value = '  -0.07%  '
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = '0.183'
$c.ClearFormats()
$ws.Range("E10").Value = '  +3.81%  '
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = '0.661'
$c.ClearFormats()
$ws.Range("E11").Value = '  +4.00%  '
$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = '58.13'
$c.ClearFormats()
$ws.Range("E12").Value = '  +8.22%  '
$ws.Range("E13").Value = '  +4.90%  '
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = '9.70'
$c.ClearFormats()
$ws.Range("E14").Value = '  +5.04%  '
$ws.Range("D15").Value = '4.152.25'
$ws.Range("E15").Value = '  +5.39%  '
$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = '19.35'
$c.ClearFormats()
$ws.Range("E16").Value = '  +5.65%  '
$ws.Range("D17").Value = '3.581.95'
$ws.Range("E17").Value = '  +5.04%  '
$ws.Range("D18").Value = '69.709.49'
$ws.Range("E18").Value = '  +5.99%  '
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = '12.58'
$c.ClearFormats()
$ws.Range("E19").Value = '  +5.91%  '
$ws.Range("E20").Value = '  +0.67%  '
$ws.Range("E21").Value = '  +4.94%  '
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = '497.93'
$c.ClearFormats()
$ws.Range("E22").Value = '  +3.54%  '
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = '5.48'
$c.ClearFormats()
$ws.Range("E23").Value = '  +11.00%  '
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = '17.14'
$c.ClearFormats()
$ws.Range("E24").Value = '  +19.35%  '
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = '4.46'
$c.ClearFormats()
$ws.Range("E25").Value = '  +8.09%  '
$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = '91.03'
$c.ClearFormats()
$ws.Range("E26").Value = '  +1.85%  '
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = '3.07'
$c.ClearFormats()
$ws.Range("E27").Value = '  +5.23%  '
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = '11.16'
$c.ClearFormats()
$ws.Range("E28").Value = '  +4.47%  '
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = '9.30'
$c.ClearFormats()
$ws.Range("E29").Value = '  +6.21%  '
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = '32.14'
$c.ClearFormats()
$ws.Range("E30").Value = '  +2.64%  '
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = '7.49'
$c.ClearFormats()
$ws.Range("E31").Value = '  +13.60%  '
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = '12.14'
$c.ClearFormats()
$ws.Range("E32").Value = '  +5.71%  '
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = '613.98'
$c.ClearFormats()
$ws.Range("E33").Value = '  +7.43%  '
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = '65.29'
$c.ClearFormats()
$ws.Range("E34").Value = '  +2.38%  '
$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = '0.115'
$c.ClearFormats()
$ws.Range("E35").Value = '  +6.63%  '
$ws.Range("D36").Value = '0.0₃0832'
$ws.Range("E36").Value = '  +11.60%  '
$ws.Range("B37").Value = 'Stacks'
$ws.Range("C37").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = '3.75'
$c.ClearFormats()
$ws.Range("E37").Value = '  +2.70%  '
$ws.Range("B38").Value = 'Kaspa'
$ws.Range("C38").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = '0.148'
$c.ClearFormats()
$ws.Range("E38").Value = '  +4.07%  '
$ws.Range("E39").Value = '  -0.02%  '
$ws.Range("B40").Value = 'InjectiveProtocol'
$ws.Range("C40").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = '38.00'
$c.ClearFormats()
$ws.Range("E40").Value = '  +5.80%  '
$ws.Range("B41").Value = 'TheGraph'
$ws.Range("C41").Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = '0.398'
$c.ClearFormats()
$ws.Range("E41").Value = '  +6.11%  '
$ws.Range("D42").Value = '3.330.04'
$ws.Range("E42").Value = '  +7.73%  '
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = '3.09'
$c.ClearFormats()
$ws.Range("E43").Value = '  +10.06%  '
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = '2.69'
$c.ClearFormats()
$ws.Range("E44").Value = '  +9.31%  '
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = '0.0442'
$c.ClearFormats()
$ws.Range("E45").Value = '  +5.95%  '
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = '2.89'
$c.ClearFormats()
$ws.Range("E46").Value = '  +17.44%  '
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = '3.31'
$c.ClearFormats()
$ws.Range("E47").Value = '  +4.17%  '
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = '0.137'
$c.ClearFormats()
$ws.Range("E48").Value = '  +2.18%  '
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = '9.12'
$c.ClearFormats()
$ws.Range("E49").Value = '  +7.67%  '
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = '3.26'
$c.ClearFormats()
$ws.Range("E50").Value = '  +5.90%  '
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = '1.00'
$c.ClearFormats()
$ws.Range("E51").Value = '  +0.27%  '
